{"js": "// POI 32 - Struggles: fix image/audio credit labels.\n//\n// 1. Remove the stray \"_GoBack\" bookmark from the \"Notes:\" paragraph\n//    near the top of the document.\n// 2. The \"(Image: C Kamana)\" paragraph that follows the text\n//    \"... Brother Neil McGurk discusses Eric Molobi and the socio-political\n//    situation within which Sacred Heart College sat in the 1980s\" is\n//    actually a photo credit, not an audio credit: relabel it to\n//    \"(Image: Wits University, Historical Papers Archive) \" and add a new\n//    \"(Audio: C Kamana)\" paragraph right after it.\n// 3. The \"(Image: C Kamana)\" paragraph that follows the \"... xenophobia.\"\n//    paragraph is really both an image AND an audio credit: relabel it to\n//    \"(Image and Audio: C Kamana)\", and re-home the \"_GoBack\" bookmark\n//    here (right after \"Audio\", before the colon).\n\n// --- Step 1: drop the old bookmark -----------------------------------\nconst goBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nif (!goBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Step 2: first \"(Image: C Kamana)\" -> photo credit + new audio line\nconst anchor1 = context.document.body.search(\n  \"Brother Neil McGurk discusses Eric Molobi\",\n  { matchCase: true }\n);\nanchor1.load(\"items\");\nawait context.sync();\nif (anchor1.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for anchor1, got \" + anchor1.items.length);\n}\n// paragraph structure: [anchor paragraph] -> [blank] -> [(Image: C Kamana)]\nlet para1 = anchor1.items[0].paragraphs.getFirst().getNext().getNext();\npara1.load(\"text\");\nawait context.sync();\nif (para1.text !== \"(Image: C Kamana)\") {\n  throw new Error(\"Unexpected paragraph text for target 1: \" + JSON.stringify(para1.text));\n}\nconst range1 = para1.getRange();\nrange1.insertText(\"(Image: Wits University, Historical Papers Archive) \", Word.InsertLocation.replace);\nawait context.sync();\npara1.insertParagraph(\"(Audio: C Kamana)\", Word.InsertLocation.after);\nawait context.sync();\n\n// --- Step 3: second \"(Image: C Kamana)\" -> image+audio credit ---------\nconst anchor2 = context.document.body.search(\"xenophobia.\", { matchCase: true });\nanchor2.load(\"items\");\nawait context.sync();\nif (anchor2.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for anchor2, got \" + anchor2.items.length);\n}\n// paragraph structure: [anchor paragraph] -> [blank] -> [(Image: C Kamana)]\nlet para2 = anchor2.items[0].paragraphs.getFirst().getNext().getNext();\npara2.load(\"text\");\nawait context.sync();\nif (para2.text !== \"(Image: C Kamana)\") {\n  throw new Error(\"Unexpected paragraph text for target 2: \" + JSON.stringify(para2.text));\n}\nconst range2 = para2.getRange();\nconst sub2 = range2.search(\"(Image:\", { matchCase: true });\nsub2.load(\"items\");\nawait context.sync();\nif (sub2.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for '(Image:' inside target 2\");\n}\nsub2.items[0].insertText(\"(Image and Audio:\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-home the bookmark right after \" and Audio\" (before the colon).\nconst range2b = para2.getRange();\nconst sub2b = range2b.search(\" and Audio\", { matchCase: true });\nsub2b.load(\"items\");\nawait context.sync();\nif (sub2b.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for ' and Audio' inside target 2\");\n}\nconst bookmarkSpot = sub2b.items[0].getRange(Word.RangeLocation.after);\nbookmarkSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# POI 32 - Struggles: fix image/audio credit labels.\n#\n# 1. Remove the stray \"_GoBack\" bookmark from the \"Notes:\" paragraph\n#    near the top of the document.\n# 2. The \"(Image: C Kamana)\" paragraph that follows the text\n#    \"... Brother Neil McGurk discusses Eric Molobi and the socio-political\n#    situation within which Sacred Heart College sat in the 1980s\" is\n#    actually a photo credit, not an audio credit: relabel it to\n#    \"(Image: Wits University, Historical Papers Archive) \" and add a new\n#    \"(Audio: C Kamana)\" paragraph right after it.\n# 3. The \"(Image: C Kamana)\" paragraph that follows the \"... xenophobia.\"\n#    paragraph is really both an image AND an audio credit: relabel it to\n#    \"(Image and Audio: C Kamana)\", and re-home the \"_GoBack\" bookmark\n#    here (right after \"Audio\", before the colon).\n\n$d = $word.ActiveDocument\n\n# --- Step 1: drop the old bookmark -----------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Step 2: first \"(Image: C Kamana)\" -> photo credit + new audio line\n$search1 = $d.Content\n$found1 = $search1.Find.Execute(\"Brother Neil McGurk discusses Eric Molobi\")\nif (-not $found1) {\n    throw \"Could not find anchor1 text\"\n}\n# paragraph structure: [anchor paragraph] -> [blank] -> [(Image: C Kamana)]\n$anchorPara1 = $search1.Paragraphs(1)\n$target1 = $anchorPara1.Next().Next()\nif ($target1.Range.Text.TrimEnd([char]13) -ne \"(Image: C Kamana)\") {\n    throw \"Unexpected paragraph text for target 1: [\" + $target1.Range.Text + \"]\"\n}\n$targetRange1 = $d.Range($target1.Range.Start, $target1.Range.End - 1)\n$targetRange1.Text = \"(Image: Wits University, Historical Papers Archive) \"\n# Insert a brand-new paragraph right after it containing the audio credit.\n$target1.Range.InsertParagraphAfter()\n$newPara1 = $target1.Next()\n$newRange1 = $d.Range($newPara1.Range.Start, $newPara1.Range.End - 1)\n$newRange1.Text = \"(Audio: C Kamana)\"\n\n# --- Step 3: second \"(Image: C Kamana)\" -> image+audio credit ---------\n$search2 = $d.Content\n$found2 = $search2.Find.Execute(\"xenophobia.\")\nif (-not $found2) {\n    throw \"Could not find anchor2 text\"\n}\n# paragraph structure: [anchor paragraph] -> [blank] -> [(Image: C Kamana)]\n$anchorPara2 = $search2.Paragraphs(1)\n$target2 = $anchorPara2.Next().Next()\nif ($target2.Range.Text.TrimEnd([char]13) -ne \"(Image: C Kamana)\") {\n    throw \"Unexpected paragraph text for target 2: [\" + $target2.Range.Text + \"]\"\n}\n$targetRange2 = $d.Range($target2.Range.Start, $target2.Range.End - 1)\n$found3 = $targetRange2.Find.Execute(\"(Image:\")\nif (-not $found3) {\n    throw \"Could not find '(Image:' inside target 2\"\n}\n$targetRange2.Text = \"(Image and Audio:\"\n\n# Re-home the bookmark right after \" and Audio\" (before the colon).\n$targetRange2b = $d.Range($target2.Range.Start, $target2.Range.End - 1)\n$found4 = $targetRange2b.Find.Execute(\" and Audio\")\nif (-not $found4) {\n    throw \"Could not find ' and Audio' inside target 2\"\n}\n$bookmarkSpot = $d.Range($targetRange2b.End, $targetRange2b.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkSpot)\n"}
